# FA_IMPORT_.xlsx edit:
# - Remove row 11 (a duplicate data row), shifting all subsequent rows up by one
#   and shrinking the used range from A1:E646 to A1:E645.
# - Update the _FilterDatabase defined name and the AutoFilter range to match.
# - Re-select the (now shifted) row 11, matching the author's post-edit selection.
# - Add conditional formatting to column B that highlights duplicate values
#   (red text on a light red fill), the classic "Duplicate Values" preset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate row (old row 11); everything below shifts up.
$ws.Rows(11).Delete()

# Re-apply AutoFilter so its stored range reflects the new extent (A1:E645).
# Toggling AutoFilterMode off first avoids Range.AutoFilter() being treated
# as "remove filter" when one is already active.
$ws.AutoFilterMode = $false
$ws.Range("A1:E645").AutoFilter() | Out-Null

# Fix up the hidden _FilterDatabase defined name to the new range.
$fdb = $wb.Names.Item("FA_IMPORT!_FilterDatabase")
$fdb.RefersTo = "=FA_IMPORT!`$A`$1:`$E`$645"

# Match the author's resulting selection: the whole of (new) row 11.
$ws.Range("A11:XFD11").Select() | Out-Null

# Highlight duplicate values in column B with the standard
# "Light Red Fill with Dark Red Text" formatting.
$rng = $ws.Range("B1:B1048576")
$fc = $rng.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615
